$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.195.11"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.904.21"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.06%  "
$st = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.59"
$ws.Range("D5").Style = $st
$ws.Range("E5").Value = "  -0.54%  "
$st = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = $st
$ws.Range("E6").Value = "  +0.01%  "
$st = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5266"
$ws.Range("D7").Style = $st
$ws.Range("E7").Value = "  +2.07%  "
$st = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3776"
$ws.Range("D8").Style = $st
$ws.Range("E8").Value = "  +1.64%  "
$st = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07254"
$ws.Range("D9").Style = $st
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E11").Value = "  -0.58%  "
$st = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08423"
$ws.Range("D12").Style = $st
$ws.Range("E12").Value = "  +10.50%  "
$ws.Range("D13").Value = "1.910.78"
$ws.Range("E13").Value = "  +0.93%  "
$st = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.74"
$ws.Range("D14").Style = $st
$ws.Range("E14").Value = "  -0.27%  "
$st = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.271"
$ws.Range("D15").Style = $st
$ws.Range("E15").Value = "  -0.13%  "
$st = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = $st
$ws.Range("E16").Value = "  +0.08%  "
$st = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008620"
$ws.Range("D17").Style = $st
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("E18").Value = "  +1.41%  "
$st = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("D19").Style = $st
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "27.225.74"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "2.151.45"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  -0.15%  "
$st = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.436"
$ws.Range("D24").Style = $st
$ws.Range("E24").Value = "  +0.07%  "
$st = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.90"
$ws.Range("D25").Style = $st
$ws.Range("E25").Value = "  +0.77%  "
$st = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.275"
$ws.Range("D26").Style = $st
$ws.Range("E26").Value = "  +5.71%  "
$st = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.750"
$ws.Range("D27").Style = $st
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("E28").Value = "  +0.52%  "
$st = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.87"
$ws.Range("D29").Style = $st
$ws.Range("E29").Value = "  +0.19%  "
$st = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.924"
$ws.Range("D30").Style = $st
$ws.Range("E30").Value = "  -1.43%  "
$st = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.814"
$ws.Range("D31").Style = $st
$ws.Range("E31").Value = "  -0.06%  "
$st = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09294"
$ws.Range("D32").Style = $st
$ws.Range("E32").Value = "  +0.93%  "
$st = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8072"
$ws.Range("D33").Style = $st
$ws.Range("E33").Value = "  +6.32%  "
$st = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05063"
$ws.Range("D34").Style = $st
$ws.Range("E34").Value = "  +0.00%  "
$st = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.237"
$ws.Range("D35").Style = $st
$ws.Range("E35").Value = "  +2.97%  "
$st = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.951"
$ws.Range("D36").Style = $st
$ws.Range("E36").Value = "  -2.64%  "
$st = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.370"
$ws.Range("D37").Style = $st
$ws.Range("E37").Value = "  +3.02%  "
$st = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.618"
$ws.Range("D38").Style = $st
$ws.Range("E38").Value = "  +1.89%  "
$st = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5726"
$ws.Range("D39").Style = $st
$ws.Range("E39").Value = "  +1.34%  "
$st = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01988"
$ws.Range("D40").Style = $st
$ws.Range("E40").Value = "  -0.45%  "
$st = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.072"
$ws.Range("D41").Style = $st
$ws.Range("E41").Value = "  -0.45%  "
$st = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.640"
$ws.Range("D42").Style = $st
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("E43").Value = "  -0.16%  "
$st = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "117.53"
$ws.Range("D44").Style = $st
$ws.Range("E44").Value = "  -0.69%  "
$st = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1515"
$ws.Range("D45").Style = $st
$ws.Range("E45").Value = "  +0.40%  "
$st = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4843"
$ws.Range("D46").Style = $st
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$st = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("D47").Style = $st
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$st = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9994"
$ws.Range("D48").Style = $st
$ws.Range("E48").Value = "  -0.03%  "
$st = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.617"
$ws.Range("D49").Style = $st
$ws.Range("E49").Value = "  +2.46%  "
$st = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.43"
$ws.Range("D50").Style = $st
$ws.Range("E50").Value = "  +0.61%  "
$st = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.67"
$ws.Range("D51").Style = $st
$ws.Range("E51").Value = "  +0.12%  "
